$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 74998.336
$ws.Range("J87").Value = 74998.336
$ws.Range("L87").Value = 74998.336
$ws.Range("N87").Value = -77494.336

$ws.Range("H88").Value = 4953.6
$ws.Range("J88").Value = 4942
$ws.Range("L88").Value = 4942
$ws.Range("N88").Value = -5754

$ws.Range("H90").Value = 74998.336
$ws.Range("J90").Value = 74998.336
$ws.Range("L90").Value = 224995.008
$ws.Range("N90").Value = -237475.008

$ws.Range("H91").Value = 4953.6
$ws.Range("J91").Value = 4942
$ws.Range("L91").Value = 4942
$ws.Range("N91").Value = -7750

$ws.Range("H138").Value = 2307.9333
$ws.Range("J138").Value = 2587.5
$ws.Range("L138").Value = 7762.5
$ws.Range("N138").Value = -18042.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 23013
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H80").Value = 59997.09
$ws.Range("J80").Value = 59997.09
$ws.Range("L80").Value = 59997.09
$ws.Range("N80").Value = -61993.09

$ws.Range("H83").Value = 59997.09
$ws.Range("J83").Value = 59997.09
$ws.Range("L83").Value = 179991.27
$ws.Range("N83").Value = -189975.27

$ws.Range("H88").Value = 1332.7368
$ws.Range("J88").Value = 1392.0769
$ws.Range("L88").Value = 1392.0769
$ws.Range("N88").Value = -2204.0769

$ws.Range("H91").Value = 1332.7368
$ws.Range("J91").Value = 1392.0769
$ws.Range("L91").Value = 1392.0769
$ws.Range("N91").Value = -4200.0769

$ws.Range("H97").Value = 76925016
$ws.Range("I97").Value = 575.6667
$ws.Range("K97").Value = 575.6667
$ws.Range("M97").Value = -79.66669999999999

$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -94060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29082.572
$ws.Range("J35").Value = 49986
$ws.Range("L35").Value = 49986
$ws.Range("N35").Value = -50606

$ws.Range("H82").Value = 15622.934
$ws.Range("J82").Value = 48897.25
$ws.Range("L82").Value = 48897.25
$ws.Range("N82").Value = -49663.25

$ws.Range("H85").Value = 15622.934
$ws.Range("J85").Value = 48897.25
$ws.Range("L85").Value = 48897.25
$ws.Range("N85").Value = -51549.25

$ws.Range("H94").Value = 91568.8
$ws.Range("I94").Value = 798
$ws.Range("J94").Value = 227725
$ws.Range("K94").Value = 798
$ws.Range("L94").Value = 227725
$ws.Range("M94").Value = -347
$ws.Range("N94").Value = -228627

$ws.Range("H99").Value = 2069.05
$ws.Range("I99").Value = 2128.6875
$ws.Range("K99").Value = 2128.6875
$ws.Range("M99").Value = -630.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 30495.084
$ws.Range("J41").Value = 43189
$ws.Range("L41").Value = 43189
$ws.Range("N41").Value = -44045

$ws.Range("H50").Value = 59967.625
$ws.Range("I50").Value = 59967
$ws.Range("J50").Value = 59967.715
$ws.Range("K50").Value = 59967
$ws.Range("L50").Value = 59967.715
$ws.Range("M50").Value = -59342
$ws.Range("N50").Value = -61217.715

$ws.Range("H51").Value = 44000
$ws.Range("J51").Value = 44000
$ws.Range("L51").Value = 44000
$ws.Range("N51").Value = -45472

$ws.Range("H58").Value = 28738482
$ws.Range("I58").Value = 20835940
$ws.Range("J58").Value = 38464690
$ws.Range("K58").Value = 20835940
$ws.Range("L58").Value = 38464690
$ws.Range("M58").Value = -20835737
$ws.Range("N58").Value = -38465096

$ws.Range("H61").Value = 44000
$ws.Range("J61").Value = 44000
$ws.Range("L61").Value = 44000
$ws.Range("N61").Value = -44696

$ws.Range("H136").Value = 28738482
$ws.Range("I136").Value = 20835940
$ws.Range("J136").Value = 38464690
$ws.Range("K136").Value = 62507820
$ws.Range("L136").Value = 115394070
$ws.Range("M136").Value = -62505270
$ws.Range("N136").Value = -115399170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25451150
$ws.Range("J4").Value = 116039.555
$ws.Range("L4").Value = 348118.665
$ws.Range("N4").Value = -348342.665

$ws.Range("H139").Value = 5079.2
$ws.Range("I139").Value = 3570.7144
$ws.Range("K139").Value = 10712.1432
$ws.Range("M139").Value = -5572.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4004737.8
$ws.Range("I22").Value = 7999
$ws.Range("J22").Value = 5003922.5
$ws.Range("K22").Value = 7999
$ws.Range("L22").Value = 5003922.5
$ws.Range("M22").Value = -7704
$ws.Range("N22").Value = -5004512.5

$ws.Range("H27").Value = 4004737.8
$ws.Range("I27").Value = 7999
$ws.Range("J27").Value = 5003922.5
$ws.Range("K27").Value = 7999
$ws.Range("L27").Value = 5003922.5
$ws.Range("M27").Value = -7892
$ws.Range("N27").Value = -5004136.5

$ws.Range("H46").Value = 1898.6552
$ws.Range("I46").Value = 734
$ws.Range("K46").Value = 734
$ws.Range("M46").Value = -546

$ws.Range("H55").Value = 2841.2354
$ws.Range("I55").Value = 675.1
$ws.Range("J55").Value = 5935.7144
$ws.Range("K55").Value = 675.1
$ws.Range("L55").Value = 5935.7144
$ws.Range("M55").Value = -502.1
$ws.Range("N55").Value = -6281.7144

$ws.Range("H93").Value = 2185.84
$ws.Range("J93").Value = 2614.1667
$ws.Range("L93").Value = 2614.1667
$ws.Range("N93").Value = -5110.1667

$ws.Range("H132").Value = 3273.5476
$ws.Range("I132").Value = 2512.0303
$ws.Range("K132").Value = 7536.090899999999
$ws.Range("M132").Value = -5006.090899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15875
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 19400
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 19400
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -20648

$ws.Range("H65").Value = 15875
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 19400
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 97000
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -103240

$ws.Range("H96").Value = 1778.3
$ws.Range("I96").Value = 1327
$ws.Range("K96").Value = 1327
$ws.Range("M96").Value = 46

$ws.Range("H122").Value = 2973.8
$ws.Range("I122").Value = 1937
$ws.Range("J122").Value = 5177
$ws.Range("K122").Value = 5811
$ws.Range("L122").Value = 15531
$ws.Range("M122").Value = -3361
$ws.Range("N122").Value = -20431
